# ClueLayout.xlsx edit: "Changed walkway condition to make walkway W instead of A"
#
# Semantics discovered from the target diff:
#   - cells that were "A"  (old walkway marker) become "W"
#   - cells that were "W"  (previously a different tile) become "O"
#   - cells that were "WD" (door variant of the old "W" tile) become "OD"
#   - every other tile code (X, R, S, C, M, U, T, J, P, MD, SU, SD, CL, RR,
#     RL, PR, PU, JR, JU, TD, UR) is left untouched.
#
# A single left-to-right, top-to-bottom pass over each cell (one read +
# one write each) is used so a cell that is rewritten from "A" to "W" is
# never re-examined and accidentally rewritten again to "O" in the same
# pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The board layout occupies A1:W22 (column X and rows 23:24 are blank
# spacer/border cells with no text).
$lastRow = 22
$lastCol = 23

for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -eq "A") {
            $cell.Value2 = "W"
        }
        elseif ($val -eq "W") {
            $cell.Value2 = "O"
        }
        elseif ($val -eq "WD") {
            $cell.Value2 = "OD"
        }
    }
}

# Reflect the author's new scroll position / selection (rows 26:31 were
# selected with A31 as the last cell reached while extending the
# selection downward, and the view had scrolled so row 17 is at the top).
[void]$ws.Range("A26:A31").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
